# "Most of tasks done" -- fill in the previously-blank Pages / BibTeX Citation
# columns for the three sources that still only had placeholder data
# (Calipino de Motul, Vocabulario K'iche', Coto Manuscript), and nudge the
# current selection the way the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sources")

# Row 3 - CM / Calipino de Motul
$ws.Range("G3").Value = "N/A"
$ws.Range("H3").Value = "@manuscript{calipino_motul,<br> title = {Calipino de Motul},<br> note = {Manuscript, page 1r}<br>}"

# Row 6 - VKO / Vocabulario K'iche' Otlateca
$ws.Range("G6").Value = "N/A"
$ws.Range("H6").Value = "@manuscript{vocabulario_kiche,<br> title = {El Vocabulario K'iche' Otlateca},<br> pages = {38r}<br>}"

# Row 7 - Coto / Coto Manuscript
$ws.Range("G7").Value = "N/A"
$ws.Range("H7").Value = '"@book{coto1983,<br>  author = {Coto},<br>  title = {Coto Manuscript},<br>  year = {1983},<br>  note = {Original manuscript from 1656},<br>  pages = {277}<br>}{coto1983,<br> author = {Coto},<br> title = {Coto Manuscript},<br> year = {1983},<br> note = {Original manuscript from 1656}<br>}"'

# Row 5 - CCM's page count cell was left-aligned like the other numeric
# "pages" cells (matches G2's style)
$ws.Range("G5").HorizontalAlignment = -4131

# Leave the cursor where the author left it before saving
$ws.Range("D11").Select() | Out-Null
